$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.2928123333333333
$ws.Cells.Item(2, 8).Value = 0.8784369999999999
$ws.Cells.Item(2, 9).Value = 0.06406943071632207
$ws.Cells.Item(2, 10).Value = 0.06406943071632207
$ws.Cells.Item(2, 13).Value = 21.33926
$ws.Cells.Item(2, 14).Value = 64.01778
$ws.Cells.Item(2, 15).Value = 0.4398914187744692
$ws.Cells.Item(2, 16).Value = 0.4398914187744692
$ws.Cells.Item(2, 17).Value = 6.248398512206665
$ws.Cells.Item(2, 18).Value = 56.23558660985999
$ws.Cells.Item(2, 19).Value = 0.02818359277787547
$ws.Cells.Item(2, 20).Value = 0.02818359277787547

$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.2928123333333333
$ws.Cells.Item(3, 8).Value = 0.8784369999999999
$ws.Cells.Item(3, 9).Value = 0.06406943071632207
$ws.Cells.Item(3, 10).Value = 0.06406943071632207
$ws.Cells.Item(3, 15).Value = 0.23906065069302
$ws.Cells.Item(3, 16).Value = 0.23906065069302
$ws.Cells.Item(3, 17).Value = 3.395715738849777
$ws.Cells.Item(3, 18).Value = 30.561441649648
$ws.Cells.Item(3, 19).Value = 0.01531647979657532
$ws.Cells.Item(3, 20).Value = 0.01531647979657532

$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.2928123333333333
$ws.Cells.Item(4, 8).Value = 0.8784369999999999
$ws.Cells.Item(4, 9).Value = 0.06406943071632207
$ws.Cells.Item(4, 10).Value = 0.06406943071632207
$ws.Cells.Item(4, 13).Value = 10.59425366666667
$ws.Cells.Item(4, 14).Value = 31.782761
$ws.Cells.Item(4, 15).Value = 0.2183918878295978
$ws.Cells.Item(4, 16).Value = 0.2183918878295978
$ws.Cells.Item(4, 17).Value = 3.102128136061888
$ws.Cells.Item(4, 18).Value = 27.919153224557
$ws.Cells.Item(4, 19).Value = 0.0139922439263052
$ws.Cells.Item(4, 20).Value = 0.0139922439263052

$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.2928123333333333
$ws.Cells.Item(5, 8).Value = 0.8784369999999999
$ws.Cells.Item(5, 9).Value = 0.06406943071632207
$ws.Cells.Item(5, 10).Value = 0.06406943071632207
$ws.Cells.Item(5, 13).Value = 4.979874333333333
$ws.Cells.Item(5, 14).Value = 14.939623
$ws.Cells.Item(5, 15).Value = 0.102656042702913
$ws.Cells.Item(5, 16).Value = 0.102656042702913
$ws.Cells.Item(5, 17).Value = 1.458168623250111
$ws.Cells.Item(5, 18).Value = 13.123517609251
$ws.Cells.Item(5, 19).Value = 0.006577114215566086
$ws.Cells.Item(5, 20).Value = 0.006577114215566086

$ws.Cells.Item(6, 9).Value = 0.8630927339690215
$ws.Cells.Item(6, 10).Value = 0.8630927339690215
$ws.Cells.Item(6, 13).Value = 21.33926
$ws.Cells.Item(6, 14).Value = 64.01778
$ws.Cells.Item(6, 15).Value = 0.4398914187744692
$ws.Cells.Item(6, 16).Value = 0.4398914187744692
$ws.Cells.Item(6, 17).Value = 84.17348639644666
$ws.Cells.Item(6, 18).Value = 757.56137756802
$ws.Cells.Item(6, 19).Value = 0.3796670872795684
$ws.Cells.Item(6, 20).Value = 0.3796670872795684

$ws.Cells.Item(7, 9).Value = 0.8630927339690215
$ws.Cells.Item(7, 10).Value = 0.8630927339690215
$ws.Cells.Item(7, 15).Value = 0.23906065069302
$ws.Cells.Item(7, 16).Value = 0.23906065069302
$ws.Cells.Item(7, 19).Value = 0.2063315105910519
$ws.Cells.Item(7, 20).Value = 0.2063315105910519

$ws.Cells.Item(8, 9).Value = 0.8630927339690215
$ws.Cells.Item(8, 10).Value = 0.8630927339690215
$ws.Cells.Item(8, 13).Value = 10.59425366666667
$ws.Cells.Item(8, 14).Value = 31.782761
$ws.Cells.Item(8, 15).Value = 0.2183918878295978
$ws.Cells.Item(8, 16).Value = 0.2183918878295978
$ws.Cells.Item(8, 17).Value = 41.78941851271656
$ws.Cells.Item(8, 18).Value = 376.104766614449
$ws.Cells.Item(8, 19).Value = 0.1884924515435034
$ws.Cells.Item(8, 20).Value = 0.1884924515435034

$ws.Cells.Item(9, 9).Value = 0.8630927339690215
$ws.Cells.Item(9, 10).Value = 0.8630927339690215
$ws.Cells.Item(9, 13).Value = 4.979874333333333
$ws.Cells.Item(9, 14).Value = 14.939623
$ws.Cells.Item(9, 15).Value = 0.102656042702913
$ws.Cells.Item(9, 16).Value = 0.102656042702913
$ws.Cells.Item(9, 17).Value = 19.64329524326744
$ws.Cells.Item(9, 18).Value = 176.789657189407
$ws.Cells.Item(9, 19).Value = 0.08860168455489784
$ws.Cells.Item(9, 20).Value = 0.08860168455489784

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.332886
$ws.Cells.Item(10, 8).Value = 0.998658
$ws.Cells.Item(10, 9).Value = 0.07283783531465635
$ws.Cells.Item(10, 10).Value = 0.07283783531465635
$ws.Cells.Item(10, 13).Value = 21.33926
$ws.Cells.Item(10, 14).Value = 64.01778
$ws.Cells.Item(10, 15).Value = 0.4398914187744692
$ws.Cells.Item(10, 16).Value = 0.4398914187744692
$ws.Cells.Item(10, 17).Value = 7.10354090436
$ws.Cells.Item(10, 18).Value = 63.93186813924
$ws.Cells.Item(10, 19).Value = 0.03204073871702531
$ws.Cells.Item(10, 20).Value = 0.03204073871702532

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.332886
$ws.Cells.Item(11, 8).Value = 0.998658
$ws.Cells.Item(11, 9).Value = 0.07283783531465635
$ws.Cells.Item(11, 10).Value = 0.07283783531465635
$ws.Cells.Item(11, 15).Value = 0.23906065069302
$ws.Cells.Item(11, 16).Value = 0.23906065069302
$ws.Cells.Item(11, 17).Value = 3.860446097248
$ws.Cells.Item(11, 18).Value = 34.744014875232
$ws.Cells.Item(11, 19).Value = 0.01741266030539278
$ws.Cells.Item(11, 20).Value = 0.01741266030539278

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.332886
$ws.Cells.Item(12, 8).Value = 0.998658
$ws.Cells.Item(12, 9).Value = 0.07283783531465635
$ws.Cells.Item(12, 10).Value = 0.07283783531465635
$ws.Cells.Item(12, 13).Value = 10.59425366666667
$ws.Cells.Item(12, 14).Value = 31.782761
$ws.Cells.Item(12, 15).Value = 0.2183918878295978
$ws.Cells.Item(12, 16).Value = 0.2183918878295978
$ws.Cells.Item(12, 17).Value = 3.526678726082
$ws.Cells.Item(12, 18).Value = 31.740108534738
$ws.Cells.Item(12, 19).Value = 0.01590719235978914
$ws.Cells.Item(12, 20).Value = 0.01590719235978914

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.332886
$ws.Cells.Item(13, 8).Value = 0.998658
$ws.Cells.Item(13, 9).Value = 0.07283783531465635
$ws.Cells.Item(13, 10).Value = 0.07283783531465635
$ws.Cells.Item(13, 13).Value = 4.979874333333333
$ws.Cells.Item(13, 14).Value = 14.939623
$ws.Cells.Item(13, 15).Value = 0.102656042702913
$ws.Cells.Item(13, 16).Value = 0.102656042702913
$ws.Cells.Item(13, 17).Value = 1.657730447326
$ws.Cells.Item(13, 18).Value = 14.919574025934
$ws.Cells.Item(13, 19).Value = 0.007477243932449109
$ws.Cells.Item(13, 20).Value = 0.007477243932449109
